$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 131; this shifts existing rows 131-244 down to 132-245,
# preserving their content/formatting (format copied from the row above by default).
$ws.Rows(131).Insert()

# Populate the newly inserted row 131 with the new data record.
$ws.Range("A131").Value() = 11
$ws.Range("B131").Value() = "Vega Monumental Concepción"
$ws.Range("C131").Value() = "Bíobío"
$ws.Range("D131").Value() = 44944
$ws.Range("E131").Value() = 8
$ws.Range("F131").Value() = 100112003
$ws.Range("G131").Value() = "Ajo"
$ws.Range("H131").Value() = "Chino"
$ws.Range("I131").Value() = "Primera"
$ws.Range("J131").Value() = 220
$ws.Range("K131").Value() = 14000
$ws.Range("L131").Value() = 15000
$ws.Range("M131").Value() = 14455
$ws.Range("N131").Value() = "$/caja 10 kilos"
$ws.Range("O131").Value() = "China"
$ws.Range("P131").Value() = 1446
$ws.Range("Q131").Value() = 10
$ws.Range("R131").Value() = "Hortaliza"
